# Add blood and urine SOCs to all years
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SSD")

# Populate new shared-string values in the same order the original
# author entered them (column D, then H, for each new row, then
# finally the B-column SOC identifiers) so the shared string table
# ends up in the same order as the canonical edit.
$ws.Range("D10").Value = "??blood"
$ws.Range("H10").Value = "Blood Sample"
$ws.Range("D11").Value = "??urine"
$ws.Range("H11").Value = "Urine Sample"
$ws.Range("B10").Value = "nhanes-kb:SOC-NHANES-2013-2014-BLOOD"
$ws.Range("B11").Value = "nhanes-kb:SOC-NHANES-2013-2014-URINE"

# Remaining columns reuse existing shared strings
$ws.Range("I10").Value = "nhanes-kb:STD-NHANES-2013-2014"
$ws.Range("J10").Value = "nhanes-kb:SOC-NHANES-2013-2014-SUBJECTS"
$ws.Range("I11").Value = "nhanes-kb:STD-NHANES-2013-2014"
$ws.Range("J11").Value = "nhanes-kb:SOC-NHANES-2013-2014-SUBJECTS"

# Update the selected cell to match the final state
$ws.Range("B7").Select()
